$wb = $excel.ActiveWorkbook

$statusText = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41d57dd1735f7e26850d37c58eea62075fe95d23/e2e/540c318b-f4bd-4661-9637-144ead758457.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff1bb995472684ff0602ccc729633995f65e77ed/e2e/540c318b-f4bd-4661-9637-144ead758457.md."
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ff1bb995472684ff0602ccc729633995f65e77ed/e2e/540c318b-f4bd-4661-9637-144ead758457.md"
$displayName = "540c318b-f4bd-4661-9637-144ead758457.md"

# zh-cn sheet: row 7 (540c318b-...) just finished handback generation
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $hyperlinkTarget, "", "", $displayName)
$wsZh.Range("J7").Value = "540c318b-f4bd-4661-9637-144ead758457.d99b7fac9eed6536143c7d078c233d53c93168ec.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-29 18:59:53"
$wsZh.Range("P7").Value = $statusText

# de-de sheet: row 7 (540c318b-...) just finished handback generation
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $hyperlinkTarget, "", "", $displayName)
$wsDe.Range("J7").Value = "540c318b-f4bd-4661-9637-144ead758457.d99b7fac9eed6536143c7d078c233d53c93168ec.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-29 19:00:04"
$wsDe.Range("P7").Value = $statusText
